# Add municipality tax data for St.Gallen (canton SG, BfsId 3203)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: 2019 data for St.Gallen --------------------------------------
$ws.Range("A4").Value = 3203
$ws.Range("B4").Value = "SG"
$ws.Range("C4").Value = 2019
$ws.Range("D4").Value = "St.Gallen"
$ws.Range("E4").Value = 115
$ws.Range("F4").Value = 141
$ws.Range("G4").Value = 25
$ws.Range("H4").Value = 25
$ws.Range("I4").Value = 26
$ws.Range("J4").Value = 0
$ws.Range("J4").NumberFormat = "0.00"

# --- Row 5: 2018 data for St.Gallen --------------------------------------
$ws.Range("A5").Value = 3203
$ws.Range("B5").Value = "SG"
$ws.Range("C5").Value = 2018
$ws.Range("D5").Value = "St.Gallen"
$ws.Range("E5").Value = 115
$ws.Range("F5").Value = 144
$ws.Range("G5").Value = 25
$ws.Range("H5").Value = 25
$ws.Range("I5").Value = 26
$ws.Range("J5").Value = 0
$ws.Range("J5").NumberFormat = "0.00"

# --- Page setup (A4 paper, portrait) --------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Restore the active cell/selection on the sheet -----------------------
$ws.Range("H11").Select() | Out-Null
